# ============================================================
# 1) "总计" (summary) sheet: insert a 2022-Q3 row above 2022-Q2
# ============================================================
$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("总计")

# Shift the existing quarters down one row, opening up row 2 for 2022-Q3.
$summary.Rows.Item(2).Insert()

# Row-insert leaves the new row with odd inherited formatting; re-sync it with
# the rest of the index column (style used by A3:A6) and drop the stray format
# that landed on B2:D2 so they go back to "no explicit style", like the others.
$summary.Range("A3").Copy($summary.Range("A2"))
$summary.Range("B2:D2").ClearFormats()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 28
$summary.Cells.Item(2, 4).Value = 2.49

# ============================================================
# 2) brand-new "2022-Q3" sheet, positioned right after "总计"
# ============================================================
$totalSheet = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# Clone the header/formatting skeleton from the structurally-identical "2022-Q2"
# sheet (bold+bordered header row, bold index column) so the new tab ends up
# with the exact same look, then fully overwrite the values below.
$template = $wb.Worksheets.Item("2022-Q2")
$template.Range("A1:H29").Copy($q3.Range("A1:H29"))
$template.Range("A2").Copy($q3.Range("A26:A29"))
$q3.Range("A1").ClearContents()

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$rows = @(
    @(0, "160813", "长盛同盛成长优选灵活配置混合（LOF）", "11.34", "86.61", "3.74", "0.4241", 7),
    @(1, "010885", "长盛优势企业精选混合A", "10.61", "84.51", "3.50", "0.3714", 9),
    @(2, "519039", "长盛同德主题混合", "9.65", "83.31", "3.52", "0.3397", 8),
    @(3, "000534", "长盛高端装备制造灵活配置混合", "3.94", "80.00", "7.64", "0.3010", 2),
    @(4, "009800", "长盛制造精选混合A", "4.65", "83.23", "3.51", "0.1632", 8),
    @(5, "000598", "长盛生态环境主题灵活配置混合", "1.75", "84.44", "7.80", "0.1365", 2),
    @(6, "014325", "国联安核心趋势一年持有混合A", "3.74", "86.69", "3.27", "0.1223", 10),
    @(7, "001892", "长盛新兴成长主题灵活配置混合", "1.33", "71.07", "7.95", "0.1057", 2),
    @(8, "010155", "长盛核心成长混合A", "2.52", "81.05", "3.53", "0.0890", 7),
    @(9, "002156", "长盛盛世灵活配置混合A", "2.29", "30.91", "2.88", "0.0660", 2),
    @(10, "080002", "长盛创新先锋混合A", "0.77", "76.75", "7.17", "0.0552", 2),
    @(11, "014885", "长盛匠心研究混合A", "1.54", "48.22", "2.67", "0.0411", 5),
    @(12, "005265", "博时厚泽回报灵活配置混合A", "1.70", "72.26", "2.35", "0.0400", 10),
    @(13, "257050", "国联安主题驱动混合", "1.48", "93.02", "2.62", "0.0388", 10),
    @(14, "004332", "恒生前海沪港深新兴产业精选混合", "0.47", "92.74", "6.95", "0.0327", 5),
    @(15, "014886", "长盛匠心研究混合C", "1.13", "48.22", "2.67", "0.0302", 5),
    @(16, "010886", "长盛优势企业精选混合C", "0.75", "84.51", "3.50", "0.0262", 9),
    @(17, "005266", "博时厚泽回报灵活配置混合C", "0.96", "72.26", "2.35", "0.0226", 10),
    @(18, "000410", "益民服务领先灵活配置混合", "0.71", "92.80", "2.91", "0.0207", 10),
    @(19, "010156", "长盛核心成长混合C", "0.45", "81.05", "3.53", "0.0159", 7),
    @(20, "002085", "长盛互联网+主题灵活配置混合", "0.31", "83.56", "5.11", "0.0158", 5),
    @(21, "014326", "国联安核心趋势一年持有混合C", "0.38", "86.69", "3.27", "0.0124", 10),
    @(22, "009801", "长盛制造精选混合C", "0.20", "83.23", "3.51", "0.0070", 8),
    @(23, "000804", "中信建投稳利混合A", "0.21", "38.27", "2.40", "0.0050", 4),
    @(24, "012716", "长盛创新先锋混合C", "0.05", "76.75", "7.17", "0.0036", 2),
    @(25, "006844", "中信建投稳利混合C", "0.10", "38.27", "2.40", "0.0024", 4),
    @(26, "006084", "融通研究优选混合", "0.07", "77.79", "2.30", "0.0016", 10),
    @(27, "002157", "长盛盛世灵活配置混合C", "0.05", "30.91", "2.88", "0.0014", 2)
)

# Pre-format the numeric-looking text columns (fund code + the size/position/ratio/
# value figures) as Text so leading zeros and trailing-zero precision survive --
# a plain .Value assignment would otherwise coerce "010885" -> 10885, "3.50" -> 3.5.
$q3.Range("B2:G29").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $q3.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}
